# aggiornamento a 9/09 compreso
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (dates 2021-09-02 .. 2021-09-09 == serials 44441..44448)
$newRows = @(
    @(44441, 22, 118, 163.0532410286172),
    @(44442, 5,  112, 154.7623982644503),
    @(44443, 14, 105, 145.0897483729221),
    @(44444, 2,  103, 142.3261341181998),
    @(44445, 22, 76,  105.0173416794484),
    @(44446, 5,  72,  99.49011317000372),
    @(44447, 4,  74,  102.2537274247261),
    @(44448, 18, 70,  96.7264989152814)
)

$lastRow = 366
$startRow = $lastRow + 1

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    # Copy the formatting from the cell directly above so the new date
    # cell picks up the same style (border, bold, centered, date format)
    # without introducing any new style entries.
    $ws.Cells.Item($lastRow, 1).Copy()
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.PasteSpecial(-4122)
    $cellA.Value = $rowData[0]

    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
}

$excel.CutCopyMode = $false
